# Add a new "Estimates" worksheet with a story-point estimate summary,
# matching the author's commit that appended this sheet to the workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create the new worksheet directly after Sheet1.
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Estimates"

# Populate the summary cells.
$label = "Total Story Points Esimates (Including Desing, Cut Effort, DB Design, Testing, Requirement Detailing, Code Review, Bug Fixing, Documentation, Release Notes)"
$newSheet.Range("B2").Value = $label
$newSheet.Range("B2").WrapText = $true
$newSheet.Range("C2").Value = 314

# Column widths as seen in the final workbook.
$newSheet.Columns.Item(2).ColumnWidth = 56
$newSheet.Columns.Item(3).ColumnWidth = 10.54296875

# Row height for the populated row.
$newSheet.Rows.Item(2).RowHeight = 43.5

# Make the new sheet the active one (it becomes the visible/selected tab).
$newSheet.Activate()
$newSheet.Range("C2").Select()
